$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.142.05"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "2.528.47"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'581.38"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "'152.07"
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "2.532.48"
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "'29.59"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "2.989.17"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "63.919.61"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "2.533.73"
$ws.Range("E18").Value = "  +3.10%  "
$ws.Range("D19").Value = "'7.85"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "'10.96"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  +2.88%  "
$ws.Range("D22").Value = "'328.09"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'10.08"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "'65.42"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "'653.84"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.49"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("D32").Value = "'8.03"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "'0.136"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").Value = "'5.53"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").Value = "'0.373"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "'18.89"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'151.98"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.80"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "'1.78"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").Value = "'162.46"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "0.0₆0303"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "'15.42"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'3.65"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").Value = "'20.98"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").Value = "'0.618"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("E51").Value = "  +0.86%  "
